$wb = $excel.ActiveWorkbook
$unitsSheet = $wb.Worksheets.Item("Units")

# Add a new worksheet right after "Units" and name it "Basics".
$ws = $wb.Worksheets.Add($null, $unitsSheet)
$ws.Name = "Basics"

# Populate the basic balance constants.
$ws.Range("A1").Value = "HITPOINTS_BASE"
$ws.Range("B1").Value = 40

$ws.Range("A2").Value = "HITPOINTS_PER_LEVEL_BASE"
$ws.Range("B2").Formula = "=B1/10"

$ws.Range("A3").Value = "ARMOR_BASE"
$ws.Range("B3").Value = 0.5

$ws.Range("A4").Value = "DAMAGE_BASE"
$ws.Range("B4").Value = 5

$ws.Range("A5").Value = "SUPPORT_BASE"
$ws.Range("B5").Value = 2

$ws.Range("A6").Value = "HERO_POWER_MULTIPLICATOR"
$ws.Range("B6").Value = 2.5

$ws.Range("A7").Value = "POWER_BASE"
$ws.Range("B7").Value = 1

$ws.Range("A8").Value = "POWER_PER_LEVEL_BASE"
$ws.Range("B8").Value = 0.1

# Widen the label column (stored OOXML width ends up as 32).
$ws.Columns.Item(1).ColumnWidth = 31.166666666666668

# Match the page setup used for the new sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Move the selection on "Units" off of its old spot, and select a cell
# on the new sheet, which also becomes the active tab.
$unitsSheet.Range("B2").Select()
$ws.Range("C6").Select()
